# Auto-generated edit script: applies the Tiamat_Profits.xlsx row updates
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 838.3200000000001
$ws.Range("I15").Value = 838.3200000000001
$ws.Range("K15").Value = 2514.96
$ws.Range("M15").Value = -2345.96
# Row 92
$ws.Range("H92").Value = 20000412
$ws.Range("I92").Value = 23809924
$ws.Range("J92").Value = 476.25
$ws.Range("K92").Value = 23809924
$ws.Range("L92").Value = 476.25
$ws.Range("M92").Value = -23808676
$ws.Range("N92").Value = -2972.25
# Row 137
$ws.Range("H137").Value = 24388.49
$ws.Range("I137").Value = 26496.744
$ws.Range("K137").Value = 79490.23199999999
$ws.Range("M137").Value = -76940.23199999999

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 23737
$ws.Range("I3").Value = 602.5
$ws.Range("J3").Value = 70006
$ws.Range("K3").Value = 602.5
$ws.Range("L3").Value = 70006
$ws.Range("M3").Value = -487.5
$ws.Range("N3").Value = -70236
# Row 6
$ws.Range("H6").Value = 44501.1
$ws.Range("I6").Value = 1002
$ws.Range("K6").Value = 1002
$ws.Range("M6").Value = -829
# Row 24
$ws.Range("H24").Value = 48405
$ws.Range("I24").Value = 45000
$ws.Range("J24").Value = 48972.5
$ws.Range("K24").Value = 45000
$ws.Range("L24").Value = 48972.5
$ws.Range("M24").Value = -44626
$ws.Range("N24").Value = -49720.5
# Row 32
$ws.Range("H32").Value = 188339.08
$ws.Range("I32").Value = 192315.53
$ws.Range("J32").Value = 158231.72
$ws.Range("K32").Value = 192315.53
$ws.Range("L32").Value = 158231.72
$ws.Range("M32").Value = -192028.53
$ws.Range("N32").Value = -158805.72
# Row 74
$ws.Range("H74").Value = 45690.42
$ws.Range("I74").Value = 81340.84
$ws.Range("K74").Value = 81340.84
$ws.Range("M74").Value = -80466.84
# Row 77
$ws.Range("H77").Value = 45690.42
$ws.Range("I77").Value = 81340.84
$ws.Range("K77").Value = 406704.2
$ws.Range("M77").Value = -402336.2
# Row 88
$ws.Range("H88").Value = 23640
$ws.Range("I88").Value = 5060
$ws.Range("J88").Value = 42220
$ws.Range("K88").Value = 5060
$ws.Range("L88").Value = 42220
$ws.Range("M88").Value = -4654
$ws.Range("N88").Value = -43032
# Row 91
$ws.Range("H91").Value = 23640
$ws.Range("I91").Value = 5060
$ws.Range("J91").Value = 42220
$ws.Range("K91").Value = 5060
$ws.Range("L91").Value = 42220
$ws.Range("M91").Value = -3656
$ws.Range("N91").Value = -45028
# Row 100
$ws.Range("H100").Value = 48405
$ws.Range("I100").Value = 45000
$ws.Range("J100").Value = 48972.5
$ws.Range("K100").Value = 45000
$ws.Range("L100").Value = 48972.5
$ws.Range("M100").Value = -43918
$ws.Range("N100").Value = -51136.5
# Row 101
$ws.Range("H101").Value = 39985
$ws.Range("J101").Value = 39985
$ws.Range("L101").Value = 39985
$ws.Range("N101").Value = -46475

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 497
$ws.Range("I11").Value = 497
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 497
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -357
$ws.Range("N11").ClearContents()
# Row 12
$ws.Range("H12").Value = 3398.3333
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5336
# Row 86
$ws.Range("H86").Value = 779611.1
$ws.Range("I86").Value = 2083.3333
$ws.Range("J86").Value = 2334666.8
$ws.Range("K86").Value = 2083.3333
$ws.Range("L86").Value = 2334666.8
$ws.Range("M86").Value = -960.3332999999998
$ws.Range("N86").Value = -2336912.8
# Row 89
$ws.Range("H89").Value = 779611.1
$ws.Range("I89").Value = 2083.3333
$ws.Range("J89").Value = 2334666.8
$ws.Range("K89").Value = 10416.6665
$ws.Range("L89").Value = 11673334
$ws.Range("M89").Value = -4800.666499999999
$ws.Range("N89").Value = -11684566
# Row 100
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
# Row 134
$ws.Range("H134").Value = 23304.318
$ws.Range("I134").Value = 1223.9166
$ws.Range("J134").Value = 95567.45
$ws.Range("K134").Value = 3671.7498
$ws.Range("L134").Value = 286702.35
$ws.Range("M134").Value = -1136.7498
$ws.Range("N134").Value = -291772.35

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 35000.668
$ws.Range("J4").Value = 35000.668
$ws.Range("L4").Value = 35000.668
$ws.Range("N4").Value = -35224.668
# Row 13
$ws.Range("H13").Value = 23500.5
$ws.Range("J13").Value = 23500.5
$ws.Range("L13").Value = 23500.5
$ws.Range("N13").Value = -23778.5
# Row 58
$ws.Range("H58").Value = 4137.579
$ws.Range("I58").Value = 1221.2
$ws.Range("J58").Value = 7378
$ws.Range("K58").Value = 1221.2
$ws.Range("L58").Value = 7378
$ws.Range("M58").Value = -1018.2
$ws.Range("N58").Value = -7784
# Row 62
$ws.Range("H62").Value = 71431144
$ws.Range("I62").Value = 83335660
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 83335660
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -83335036
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 71431144
$ws.Range("I65").Value = 83335660
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 416678300
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -416675180
$ws.Range("N65").Value = -26240
# Row 132
$ws.Range("H132").Value = 2106.2144
$ws.Range("I132").Value = 1687.04
$ws.Range("J132").Value = 5599.3335
$ws.Range("K132").Value = 5061.12
$ws.Range("L132").Value = 16798.0005
$ws.Range("M132").Value = -2531.12
$ws.Range("N132").Value = -21858.0005
# Row 134
$ws.Range("H134").Value = 13890653
$ws.Range("I134").Value = 1418.25
$ws.Range("J134").Value = 31252196
$ws.Range("K134").Value = 4254.75
$ws.Range("L134").Value = 93756588
$ws.Range("M134").Value = -1719.75
$ws.Range("N134").Value = -93761658
# Row 136
$ws.Range("H136").Value = 4137.579
$ws.Range("I136").Value = 1221.2
$ws.Range("J136").Value = 7378
$ws.Range("K136").Value = 3663.6
$ws.Range("L136").Value = 22134
$ws.Range("M136").Value = -1113.6
$ws.Range("N136").Value = -27234

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
# Row 13
$ws.Range("H13").Value = 200
$ws.Range("I13").Value = 200
$ws.Range("K13").Value = 600
$ws.Range("M13").Value = -432
# Row 23
$ws.Range("H23").Value = 68
$ws.Range("I23").Value = 49
$ws.Range("J23").Value = 80.666664
$ws.Range("K23").Value = 147
$ws.Range("L23").Value = 241.999992
$ws.Range("M23").Value = 88
$ws.Range("N23").Value = -711.999992
# Row 32
$ws.Range("H32").Value = 1847.4736
$ws.Range("J32").Value = 1394.4445
$ws.Range("L32").Value = 4183.333500000001
$ws.Range("N32").Value = -4749.333500000001
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 22666.666
$ws.Range("J4").Value = 22666.666
$ws.Range("L4").Value = 22666.666
$ws.Range("N4").Value = -22890.666
# Row 5
$ws.Range("H5").Value = 12000.5
$ws.Range("I5").Value = 8004
$ws.Range("J5").Value = 12571.429
$ws.Range("K5").Value = 8004
$ws.Range("L5").Value = 12571.429
$ws.Range("M5").Value = -7892
$ws.Range("N5").Value = -12795.429
# Row 11
$ws.Range("H11").Value = 26669112
$ws.Range("I11").Value = 48000400
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 48000400
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -48000261
$ws.Range("N11").Value = -5278
# Row 132
$ws.Range("H132").Value = 54500.844
$ws.Range("I132").Value = 1466.9333
$ws.Range("K132").Value = 4400.7999
$ws.Range("M132").Value = -1870.7999

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 355545.88
$ws.Range("I132").Value = 85637.75
$ws.Range("J132").Value = 1003325.4
$ws.Range("K132").Value = 256913.25
$ws.Range("L132").Value = 3009976.2
$ws.Range("M132").Value = -254383.25
$ws.Range("N132").Value = -3015036.2

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4060.718
$ws.Range("I122").Value = 3449.111
$ws.Range("K122").Value = 10347.333
$ws.Range("M122").Value = -7897.332999999999
# Row 132
$ws.Range("H132").Value = 7811
$ws.Range("I132").Value = 1787.5
$ws.Range("K132").Value = 5362.5
$ws.Range("M132").Value = -2832.5
